$d = $word.ActiveDocument

# The target run reads ": … " right after "... biểu quyết hợp lệ" and
# right before "phiếu" (proofErr-wrapped). It is the only paragraph in the
# document where that exact preceding phrase ("biểu quyết hợp lệ") is
# followed by the "…" placeholder, so anchor the Find on that unique text
# to be robust against any paragraph re-numbering.
$searchText = "biểu quyết hợp lệ: … "

$range = $d.Content
$found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target placeholder text"
}

# $range now covers "... hợp lệ: … " (Find leaves the match selected).
# Narrow it down to just the ": … " tail (the run we need to rewrite) so we
# keep the preceding label run untouched and only replace this run's text.
$tail = $range.Duplicate
[void]$tail.MoveStart(1, ($range.Text.Length - 4))

# Replace the text in place (Range.Text, not Find.Replace, so AutoCorrect's
# smart-quote substitution never touches the literal single quotes below).
$tail.Text = ": {pending_approve_total_capital | divideBy: '10000'} "
